# Updates the "cryptos" price list on Sheet1 with the latest scrape values.
# For each changed row we refresh the Price (column D) and Volume(1h)
# (column E) text, and for rows 47-48 the coin swaps rank (Monero now
# above Maker), so Coin (B), Link (C), Price (D) and Volume (E) all move.
#
# The Price column holds numbers formatted as plain text (e.g. "446.00",
# "65.646.71"), so for values that Excel would otherwise auto-parse as a
# real number (and silently drop the trailing zero / reformat), we briefly
# force a text number format before writing the value, then restore the
# cell's normal style so the sheet's styling stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range('D2').Value = '65.646.71'
    $ws.Range('E2').Value = '  -0.24%  '
    # Row 3
    $ws.Range('D3').Value = '2.944.24'
    $ws.Range('E3').Value = '  -2.33%  '
    # Row 4
    $ws.Range('E4').Value = '  +0.07%  '
    # Row 5
    $ws.Range('D5').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D5').Value = '570.69'
    $ws.Range('D5').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E5').Value = '  -2.27%  '
    # Row 6
    $ws.Range('D6').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D6').Value = '162.91'
    $ws.Range('D6').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E6').Value = '  +0.74%  '
    # Row 8
    $ws.Range('E8').Value = '  -0.54%  '
    # Row 9
    $ws.Range('D9').Value = '2.939.75'
    $ws.Range('E9').Value = '  -2.44%  '
    # Row 10
    $ws.Range('D10').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D10').Value = '6.70'
    $ws.Range('D10').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E10').Value = '  -1.22%  '
    # Row 11
    $ws.Range('E11').Value = '  -4.05%  '
    # Row 12
    $ws.Range('E12').Value = '  +0.76%  '
    # Row 13
    $ws.Range('E13').Value = '  -3.85%  '
    # Row 14
    $ws.Range('D14').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D14').Value = '34.98'
    $ws.Range('D14').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E14').Value = '  +0.66%  '
    # Row 15
    $ws.Range('E15').Value = '  -0.53%  '
    # Row 16
    $ws.Range('D16').Value = '65.635.72'
    $ws.Range('E16').Value = '  -0.20%  '
    # Row 17
    $ws.Range('D17').Value = '3.434.16'
    $ws.Range('E17').Value = '  -2.19%  '
    # Row 18
    $ws.Range('E18').Value = '  +1.28%  '
    # Row 19
    $ws.Range('D19').Value = '2.944.02'
    $ws.Range('E19').Value = '  -2.19%  '
    # Row 20
    $ws.Range('D20').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D20').Value = '15.81'
    $ws.Range('D20').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E20').Value = '  +12.69%  '
    # Row 21
    $ws.Range('D21').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D21').Value = '446.00'
    $ws.Range('D21').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E21').Value = '  -2.67%  '
    # Row 22
    $ws.Range('D22').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D22').Value = '0.696'
    $ws.Range('D22').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E22').Value = '  +0.90%  '
    # Row 23
    $ws.Range('D23').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D23').Value = '7.26'
    $ws.Range('D23').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E23').Value = '  -1.66%  '
    # Row 24
    $ws.Range('D24').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D24').Value = '82.04'
    $ws.Range('D24').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E24').Value = '  -0.48%  '
    # Row 25
    $ws.Range('D25').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D25').Value = '2.24'
    $ws.Range('D25').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E25').Value = '  -1.34%  '
    # Row 26
    $ws.Range('D26').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D26').Value = '12.20'
    $ws.Range('D26').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E26').Value = '  -1.35%  '
    # Row 27
    $ws.Range('E27').Value = '  -0.04%  '
    # Row 28
    $ws.Range('D28').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D28').Value = '9.98'
    $ws.Range('D28').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E28').Value = '  -6.26%  '
    # Row 29
    $ws.Range('D29').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D29').Value = '8.16'
    $ws.Range('D29').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E29').Value = '  +1.33%  '
    # Row 30
    $ws.Range('E30').Value = '  +4.82%  '
    # Row 31
    $ws.Range('E31').Value = '  -0.63%  '
    # Row 32
    $ws.Range('E32').Value = '  -6.31%  '
    # Row 33
    $ws.Range('D33').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D33').Value = '0.116'
    $ws.Range('D33').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E33').Value = '  +4.85%  '
    # Row 34
    $ws.Range('D34').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D34').Value = '27.27'
    $ws.Range('D34').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E34').Value = '  +0.56%  '
    # Row 35
    $ws.Range('D35').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D35').Value = '0.999'
    $ws.Range('D35').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E35').Value = '  +0.07%  '
    # Row 36
    $ws.Range('D36').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D36').Value = '0.968'
    $ws.Range('D36').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E36').Value = '  -2.50%  '
    # Row 37
    $ws.Range('D37').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D37').Value = '5.70'
    $ws.Range('D37').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E37').Value = '  -2.25%  '
    # Row 38
    $ws.Range('D38').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D38').Value = '46.55'
    $ws.Range('D38').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E38').Value = '  +7.25%  '
    # Row 39
    $ws.Range('D39').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D39').Value = '49.14'
    $ws.Range('D39').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E39').Value = '  -1.19%  '
    # Row 40
    $ws.Range('D40').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D40').Value = '1.98'
    $ws.Range('D40').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E40').Value = '  -8.52%  '
    # Row 41
    $ws.Range('E41').Value = '  -3.08%  '
    # Row 42
    $ws.Range('E42').Value = '  -1.29%  '
    # Row 43
    $ws.Range('E43').Value = '  -6.10%  '
    # Row 44
    $ws.Range('D44').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D44').Value = '8.47'
    $ws.Range('D44').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E44').Value = '  +0.09%  '
    # Row 45
    $ws.Range('D45').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D45').Value = '381.85'
    $ws.Range('D45').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E45').Value = '  -2.81%  '
    # Row 46
    $ws.Range('E46').Value = '  -1.24%  '
    # Row 47
    $ws.Range('B47').Value = 'Monero'
    $ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    $ws.Range('D47').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D47').Value = '134.12'
    $ws.Range('D47').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E47').Value = '  -0.38%  '
    # Row 48
    $ws.Range('B48').Value = 'Maker'
    $ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    $ws.Range('D48').Value = '2.668.78'
    $ws.Range('E48').Value = '  -4.66%  '
    # Row 49
    $ws.Range('E49').Value = '  +0.02%  '
    # Row 50
    $ws.Range('D50').NumberFormat = '@'   # keep as text, not a number
    $ws.Range('D50').Value = '23.93'
    $ws.Range('D50').Style = 'Normal'     # drop the temporary text style again
    $ws.Range('E50').Value = '  +0.91%  '
    # Row 51
    $ws.Range('E51').Value = '  +1.00%  '
